$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 612
$ws.Range("F40").Value = 0.0481

$ws.Range("B40").Value = "Employees Earning More Than Their Managers"
$ws.Hyperlinks.Add($ws.Range("B40"), "https://leetcode.com/problems/employees-earning-more-than-their-managers/")

$ws.Range("I40").Value = "https://leetcode.com/problems/employees-earning-more-than-their-managers/submissions/"
